$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Regime Atual)
$ws.Range("D2").Value = 30332.8980555797
$ws.Range("F2").Value = 28842.70558552258
$ws.Range("H2").Value = 20.35502035144294
$ws.Range("I2").Value = 416326.5901685993
$ws.Range("J2").Value = 279.3777304166899

# Row 3 (Nova Proposta)
$ws.Range("D3").Value = 30332.8980555797
$ws.Range("F3").Value = 28842.70558552258
$ws.Range("H3").Value = 20.35502035144294
$ws.Range("I3").Value = 401008.5384240549
$ws.Range("J3").Value = 269.098486592597

# Row 4 (Nova c/ Aliq. Máxima)
$ws.Range("D4").Value = 30259.6043376406
$ws.Range("F4").Value = 28769.41186758348
$ws.Range("H4").Value = 20.30583629004697
$ws.Range("I4").Value = 398275.8910160765
$ws.Range("J4").Value = 267.2647319180255
